$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 323557
$ws.Cells.Item(2, 4).Value = 412295059
$ws.Cells.Item(3, 3).Value = 260
$ws.Cells.Item(3, 4).Value = 311479
$ws.Cells.Item(4, 3).Value = 322
$ws.Cells.Item(4, 4).Value = 461207
$ws.Cells.Item(8, 3).Value = 868
$ws.Cells.Item(8, 4).Value = 1276795
$ws.Cells.Item(10, 3).Value = 117739
$ws.Cells.Item(10, 4).Value = 172518675
$ws.Cells.Item(12, 3).Value = 59972
$ws.Cells.Item(12, 4).Value = 86560844
$ws.Cells.Item(16, 3).Value = 4020
$ws.Cells.Item(16, 4).Value = 5704697
$ws.Cells.Item(20, 3).Value = 6799
$ws.Cells.Item(20, 4).Value = 9488900
$ws.Cells.Item(22, 3).Value = 78099
$ws.Cells.Item(22, 4).Value = 97341572
$ws.Cells.Item(27, 3).Value = 290
$ws.Cells.Item(27, 4).Value = 416147
$ws.Cells.Item(28, 3).Value = 32646
$ws.Cells.Item(28, 4).Value = 47783310
$ws.Cells.Item(30, 3).Value = 11562
$ws.Cells.Item(30, 4).Value = 16630812
$ws.Cells.Item(33, 3).Value = 1563
$ws.Cells.Item(33, 4).Value = 2195781
$ws.Cells.Item(35, 3).Value = 1868
$ws.Cells.Item(35, 4).Value = 2637352
$ws.Cells.Item(36, 3).Value = 97920
$ws.Cells.Item(36, 4).Value = 123186688
$ws.Cells.Item(44, 3).Value = 44579
$ws.Cells.Item(44, 4).Value = 65333498
$ws.Cells.Item(46, 3).Value = 9216
$ws.Cells.Item(46, 4).Value = 13223944
$ws.Cells.Item(48, 3).Value = 1410
$ws.Cells.Item(48, 4).Value = 1957603
$ws.Cells.Item(51, 3).Value = 2356
$ws.Cells.Item(51, 4).Value = 3291669
$ws.Cells.Item(52, 3).Value = 69580
$ws.Cells.Item(52, 4).Value = 87270642
$ws.Cells.Item(54, 3).Value = 45
$ws.Cells.Item(54, 4).Value = 62197
$ws.Cells.Item(59, 3).Value = 28354
$ws.Cells.Item(59, 4).Value = 41582953
$ws.Cells.Item(62, 3).Value = 11218
$ws.Cells.Item(62, 4).Value = 16221284
$ws.Cells.Item(64, 3).Value = 1364
$ws.Cells.Item(64, 4).Value = 1906237
$ws.Cells.Item(68, 3).Value = 1507
$ws.Cells.Item(68, 4).Value = 2112344
$ws.Cells.Item(70, 3).Value = 20657
$ws.Cells.Item(70, 4).Value = 27057267
$ws.Cells.Item(71, 3).Value = 33
$ws.Cells.Item(71, 4).Value = 43781
$ws.Cells.Item(74, 3).Value = 7640
$ws.Cells.Item(74, 4).Value = 11187874
$ws.Cells.Item(76, 3).Value = 5161
$ws.Cells.Item(76, 4).Value = 7493563
$ws.Cells.Item(78, 3).Value = 283
$ws.Cells.Item(78, 4).Value = 397583
$ws.Cells.Item(79, 3).Value = 141957
$ws.Cells.Item(79, 4).Value = 176959291
$ws.Cells.Item(85, 3).Value = 63878
$ws.Cells.Item(85, 4).Value = 93622032
$ws.Cells.Item(88, 3).Value = 29961
$ws.Cells.Item(88, 4).Value = 43341861
$ws.Cells.Item(90, 3).Value = 2741
$ws.Cells.Item(90, 4).Value = 3945957
$ws.Cells.Item(91, 3).Value = 2883
$ws.Cells.Item(91, 4).Value = 4075360
$ws.Cells.Item(92, 3).Value = 33593
$ws.Cells.Item(92, 4).Value = 45527970
$ws.Cells.Item(95, 3).Value = 30
$ws.Cells.Item(95, 4).Value = 43314
$ws.Cells.Item(96, 3).Value = 8135
$ws.Cells.Item(96, 4).Value = 11960834
$ws.Cells.Item(98, 3).Value = 7485
$ws.Cells.Item(98, 4).Value = 10862221
$ws.Cells.Item(100, 3).Value = 539
$ws.Cells.Item(100, 4).Value = 765151
$ws.Cells.Item(101, 3).Value = 500
$ws.Cells.Item(101, 4).Value = 721391
$ws.Cells.Item(102, 3).Value = 10343
$ws.Cells.Item(102, 4).Value = 15715707
$ws.Cells.Item(104, 3).Value = 2558
$ws.Cells.Item(104, 4).Value = 4157088
$ws.Cells.Item(106, 3).Value = 3455
$ws.Cells.Item(106, 4).Value = 5607959
$ws.Cells.Item(108, 3).Value = 160
$ws.Cells.Item(108, 4).Value = 262445
$ws.Cells.Item(109, 3).Value = 200
$ws.Cells.Item(109, 4).Value = 308030
$ws.Cells.Item(110, 3).Value = 142655
$ws.Cells.Item(110, 4).Value = 176406032
$ws.Cells.Item(116, 3).Value = 53045
$ws.Cells.Item(116, 4).Value = 77751453
$ws.Cells.Item(118, 3).Value = 27366
$ws.Cells.Item(118, 4).Value = 39650941
$ws.Cells.Item(119, 3).Value = 1315
$ws.Cells.Item(119, 4).Value = 1797981
$ws.Cells.Item(122, 3).Value = 2300
$ws.Cells.Item(122, 4).Value = 3231668
$ws.Cells.Item(124, 3).Value = 520045
$ws.Cells.Item(124, 4).Value = 686874241
$ws.Cells.Item(125, 3).Value = 92
$ws.Cells.Item(125, 4).Value = 122289
$ws.Cells.Item(126, 3).Value = 216
$ws.Cells.Item(126, 4).Value = 318009
$ws.Cells.Item(129, 3).Value = 1385
$ws.Cells.Item(129, 4).Value = 2052682
$ws.Cells.Item(131, 3).Value = 210034
$ws.Cells.Item(131, 4).Value = 308763592
$ws.Cells.Item(132, 3).Value = 407
$ws.Cells.Item(132, 4).Value = 607250
$ws.Cells.Item(134, 3).Value = 185500
$ws.Cells.Item(134, 4).Value = 269759855
$ws.Cells.Item(137, 3).Value = 2858
$ws.Cells.Item(137, 4).Value = 4014969
$ws.Cells.Item(139, 3).Value = 6487
$ws.Cells.Item(139, 4).Value = 9165314
$ws.Cells.Item(140, 3).Value = 4
$ws.Cells.Item(140, 4).Value = 6000
$ws.Cells.Item(142, 3).Value = 45125
$ws.Cells.Item(142, 4).Value = 60248831
$ws.Cells.Item(148, 3).Value = 14180
$ws.Cells.Item(148, 4).Value = 20790133
$ws.Cells.Item(149, 3).Value = 3800
$ws.Cells.Item(149, 4).Value = 5480002
$ws.Cells.Item(152, 3).Value = 403
$ws.Cells.Item(152, 4).Value = 579716
$ws.Cells.Item(154, 3).Value = 393
$ws.Cells.Item(154, 4).Value = 555263
$ws.Cells.Item(155, 3).Value = 17798
$ws.Cells.Item(155, 4).Value = 23523618
$ws.Cells.Item(159, 3).Value = 7255
$ws.Cells.Item(159, 4).Value = 10555975
$ws.Cells.Item(161, 3).Value = 5065
$ws.Cells.Item(161, 4).Value = 7289996
$ws.Cells.Item(164, 3).Value = 271
$ws.Cells.Item(164, 4).Value = 387864
$ws.Cells.Item(166, 3).Value = 18440
$ws.Cells.Item(166, 4).Value = 29888219
$ws.Cells.Item(167, 3).Value = 2013
$ws.Cells.Item(167, 4).Value = 3283349
$ws.Cells.Item(168, 3).Value = 275
$ws.Cells.Item(168, 4).Value = 448101
$ws.Cells.Item(170, 3).Value = 64
$ws.Cells.Item(170, 4).Value = 110690
$ws.Cells.Item(172, 3).Value = 88486
$ws.Cells.Item(172, 4).Value = 110603462
$ws.Cells.Item(177, 3).Value = 644
$ws.Cells.Item(177, 4).Value = 949348
$ws.Cells.Item(179, 3).Value = 34081
$ws.Cells.Item(179, 4).Value = 49977954
$ws.Cells.Item(181, 3).Value = 13149
$ws.Cells.Item(181, 4).Value = 18998582
$ws.Cells.Item(183, 3).Value = 1250
$ws.Cells.Item(183, 4).Value = 1749929
$ws.Cells.Item(185, 3).Value = 1685
$ws.Cells.Item(185, 4).Value = 2368656
$ws.Cells.Item(187, 3).Value = 240306
$ws.Cells.Item(187, 4).Value = 298655698
$ws.Cells.Item(195, 3).Value = 87077
$ws.Cells.Item(195, 4).Value = 127641359
$ws.Cells.Item(196, 3).Value = 96
$ws.Cells.Item(196, 4).Value = 139627
$ws.Cells.Item(198, 3).Value = 33297
$ws.Cells.Item(198, 4).Value = 47930691
$ws.Cells.Item(201, 3).Value = 5140
$ws.Cells.Item(201, 4).Value = 7318777
$ws.Cells.Item(204, 3).Value = 4977
$ws.Cells.Item(204, 4).Value = 6891752
$ws.Cells.Item(207, 3).Value = 266156
$ws.Cells.Item(207, 4).Value = 329359127
$ws.Cells.Item(214, 3).Value = 619
$ws.Cells.Item(214, 4).Value = 901378
$ws.Cells.Item(216, 3).Value = 95715
$ws.Cells.Item(216, 4).Value = 140026640
$ws.Cells.Item(219, 3).Value = 51899
$ws.Cells.Item(219, 4).Value = 75010453
$ws.Cells.Item(222, 3).Value = 4695
$ws.Cells.Item(222, 4).Value = 6592244
$ws.Cells.Item(225, 3).Value = 5887
$ws.Cells.Item(225, 4).Value = 8146147
$ws.Cells.Item(228, 3).Value = 107233
$ws.Cells.Item(228, 4).Value = 134065891
$ws.Cells.Item(229, 3).Value = 76
$ws.Cells.Item(229, 4).Value = 80550
$ws.Cells.Item(233, 3).Value = 567
$ws.Cells.Item(233, 4).Value = 828439
$ws.Cells.Item(235, 3).Value = 49818
$ws.Cells.Item(235, 4).Value = 72980383
$ws.Cells.Item(236, 3).Value = 35
$ws.Cells.Item(236, 4).Value = 50211
$ws.Cells.Item(237, 3).Value = 12554
$ws.Cells.Item(237, 4).Value = 18051518
$ws.Cells.Item(241, 3).Value = 2564
$ws.Cells.Item(241, 4).Value = 3586686
$ws.Cells.Item(242, 3).Value = 260123
$ws.Cells.Item(242, 4).Value = 328463035
$ws.Cells.Item(247, 3).Value = 12
$ws.Cells.Item(247, 4).Value = 17480
$ws.Cells.Item(250, 3).Value = 96443
$ws.Cells.Item(250, 4).Value = 141312865
$ws.Cells.Item(253, 3).Value = 65731
$ws.Cells.Item(253, 4).Value = 95259938
$ws.Cells.Item(255, 3).Value = 2422
$ws.Cells.Item(255, 4).Value = 3414819
$ws.Cells.Item(258, 3).Value = 4697
$ws.Cells.Item(258, 4).Value = 6600605
